$d = $word.ActiveDocument

# Replace all whole-word occurrences of "Tantalus" with "newTantalus".
# Word's Find/Replace (ReplaceAll) preserves run formatting correctly for
# these (plain runs / runs not wrapped in a hyperlink element).
$d.Content.Find.Execute("Tantalus", $true, $false, $false, $false, $false,
                         $true, 1, $false, "newTantalus", 2)

# "Sisyphus" lives inside a <w:hyperlink> run (styled + bold + colored).
# Replacing its Range.Text directly on the very Range object returned by
# Find tends to lose that run's direct character formatting, so: locate
# the match first, then re-open a *fresh* Range over the same span to read
# /restore the formatting explicitly around the text swap.
$searchRng = $d.Content
$found = $searchRng.Find.Execute("Sisyphus", $true, $false, $false, $false,
                                  $false, $true, 1, $false, "", 0)
if ($found) {
    $s = $searchRng.Start
    $e = $searchRng.End
    $target = $d.Range($s, $e)

    $styleName = $target.CharacterStyle.NameLocal
    $bold = $target.Bold
    $color = $target.Font.Color

    $target.Text = "newSisyphus"
    $target.Style = $styleName
    $target.Bold = $bold
    $target.Font.Color = $color
}
